# Weekly data refresh: insert this week's new observation at the top of
# the data block (row 74) and push the existing rows (74-152) down by one,
# extending the table to row 153. The worksheet holds a constant set of
# columns (A,B,C,E,F,G,I,Q,R) shared by every data row, so only the new
# row's date / volume / price / unit columns need to be populated — the
# constant columns and the pushed-down rows are carried along by the
# row insert itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 74..152 down to 75..153, extending dimension to A1:R153.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with this week's observation.
$ws.Range("A74").Value2 = 9
$ws.Range("B74").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C74").Value = "Metropolitana"
$ws.Range("D74").Value2 = 44966
$ws.Range("E74").Value2 = 13
$ws.Range("F74").Value2 = 100112022
$ws.Range("G74").Value = "Arveja Verde"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value2 = 16
$ws.Range("K74").Value2 = 27000
$ws.Range("L74").Value2 = 29000
$ws.Range("M74").Value2 = 28000
$ws.Range("N74").Value = "$/saco 25 kilos"
$ws.Range("O74").Value = "Carahue"
$ws.Range("P74").Value2 = 1120
$ws.Range("Q74").Value2 = 25
$ws.Range("R74").Value = "Hortaliza"
